$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph (paragraph 2), which sat
#    right under the title heading.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Replace the final paragraph (the old image-prompt paragraph) with
#    two paragraphs:
#      - a new bold "Play CyberCatz Free: Intergalactic Slot Game
#        Review" paragraph
#      - an italic paragraph with the (former meta-description) text
#    InsertXML on the document's *last* paragraph leaves a stray empty
#    trailing paragraph behind (Word always keeps a terminating mark),
#    so we first push a temporary paragraph after it, do the XML
#    replace on the (now not-last) target paragraph, and then delete
#    the temporary trailing paragraph mark.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastPara.Range.InsertParagraphAfter()

$target = $d.Paragraphs.Item($n)
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play CyberCatz Free: Intergalactic Slot Game Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the features of CyberCatz in this slot game review. Play for free and enjoy the unique graphics and free spins feature for a chance to win big.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($xmlFrag)

# Drop the now-redundant temporary trailing paragraph mark.
$secondToLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$cleanupRange = $d.Range($secondToLast.Range.End - 1, $d.Content.End)
$cleanupRange.Delete()
